$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (27) had the "latest" date style; now that row 28
# is the newest row, row 27 reverts to the regular datetime number format.
$ws.Range("A27").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 28.
$ws.Range("A28").Value = 45612
$ws.Range("A28").NumberFormat = "YYYY-MM-DD"
$ws.Range("B28").Value = 68
$ws.Range("C28").Value = 58
$ws.Range("D28").Value = 67
